# Insert a new data row at row 40 (pushing the existing rows 40-89 down to
# 41-90) and populate the new row with a fresh weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 40:89 down to 41:90, leaving a blank row 40 behind.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new record.
$ws.Cells.Item(40, 1).Value  = 9
$ws.Cells.Item(40, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(40, 3).Value  = "Metropolitana"
$ws.Cells.Item(40, 4).Value  = 44539
$ws.Cells.Item(40, 5).Value  = 13
$ws.Cells.Item(40, 6).Value  = 100112022
$ws.Cells.Item(40, 7).Value  = "Arveja Verde"
$ws.Cells.Item(40, 8).Value  = "Sin especificar"
$ws.Cells.Item(40, 9).Value  = "Primera"
$ws.Cells.Item(40, 10).Value = 52
$ws.Cells.Item(40, 11).Value = 16000
$ws.Cells.Item(40, 12).Value = 17000
$ws.Cells.Item(40, 13).Value = 16500
$ws.Cells.Item(40, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(40, 15).Value = "Región del Maule"
$ws.Cells.Item(40, 16).Value = 660
$ws.Cells.Item(40, 17).Value = 25
$ws.Cells.Item(40, 18).Value = "Hortaliza"
